$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new columns before column D, shifting existing D:H (and their formatting) to I:M
$ws.Range("D1:H1").EntireColumn.Insert()

# Match target column widths for the newly inserted D:H columns (28,28,29,28,28)
# (ColumnWidth setter adds a constant ~0.8333 padding vs. the stored OOXML width,
#  so subtract it to land exactly on the desired width.)
$ws.Columns.Item(4).ColumnWidth = 27.16666666666667
$ws.Columns.Item(5).ColumnWidth = 27.16666666666667
$ws.Columns.Item(6).ColumnWidth = 28.16666666666667
$ws.Columns.Item(7).ColumnWidth = 27.16666666666667
$ws.Columns.Item(8).ColumnWidth = 27.16666666666667

# Row 8: new period-label headers
$ws.Cells.Item(8, 4).Value = "6 ماهه منتهی به 1399/06"
$ws.Cells.Item(8, 5).Value = "9 ماهه منتهی به 1399/09"
$ws.Cells.Item(8, 6).Value = "12 ماهه منتهی به 1399/12"
$ws.Cells.Item(8, 7).Value = "3 ماهه منتهی به 1400/03"
$ws.Cells.Item(8, 8).Value = "6 ماهه منتهی به 1400/06"

# Row 9: new publish-date headers
$ws.Cells.Item(9, 4).Value = "1400-09-30 (4)"
$ws.Cells.Item(9, 5).Value = "1400-10-30 (2)"
$ws.Cells.Item(9, 6).Value = "1401-04-15 (8)"
$ws.Cells.Item(9, 7).Value = "1401-04-29 (2)"
$ws.Cells.Item(9, 8).Value = "1401-09-14 (4)"

# Data rows 11-27: new historical figures for columns D:H
# Row 11
$ws.Cells.Item(11, 4).Value = 43286898
$ws.Cells.Item(11, 5).Value = 83110574
$ws.Cells.Item(11, 6).Value = 133177288
$ws.Cells.Item(11, 7).Value = 49950179
$ws.Cells.Item(11, 8).Value = 76819986

# Row 12
$ws.Cells.Item(12, 4).Value = -27416933
$ws.Cells.Item(12, 5).Value = -53755646
$ws.Cells.Item(12, 6).Value = -94867910
$ws.Cells.Item(12, 7).Value = -40185447
$ws.Cells.Item(12, 8).Value = -60799934

# Row 13
$ws.Cells.Item(13, 4).Value = 15869965
$ws.Cells.Item(13, 5).Value = 29354928
$ws.Cells.Item(13, 6).Value = 38309378
$ws.Cells.Item(13, 7).Value = 9764732
$ws.Cells.Item(13, 8).Value = 16020052

# Row 14
$ws.Cells.Item(14, 4).Value = -651402
$ws.Cells.Item(14, 5).Value = -1212099
$ws.Cells.Item(14, 6).Value = -2082002
$ws.Cells.Item(14, 7).Value = -747548
$ws.Cells.Item(14, 8).Value = -1272200

# Row 15
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(15, 5).Value = 0
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(15, 7).Value = 0
$ws.Cells.Item(15, 8).Value = 0

# Row 16
$ws.Cells.Item(16, 4).Value = -740642
$ws.Cells.Item(16, 5).Value = 625177
$ws.Cells.Item(16, 6).Value = 640075
$ws.Cells.Item(16, 7).Value = -424567
$ws.Cells.Item(16, 8).Value = -1052701

# Row 17
$ws.Cells.Item(17, 4).Value = 14477921
$ws.Cells.Item(17, 5).Value = 28768006
$ws.Cells.Item(17, 6).Value = 36867451
$ws.Cells.Item(17, 7).Value = 8592617
$ws.Cells.Item(17, 8).Value = 13695151

# Row 18
$ws.Cells.Item(18, 4).Value = -1443027
$ws.Cells.Item(18, 5).Value = -2754218
$ws.Cells.Item(18, 6).Value = -3434959
$ws.Cells.Item(18, 7).Value = -1022293
$ws.Cells.Item(18, 8).Value = -1631493

# Row 19
$ws.Cells.Item(19, 4).Value = 102923
$ws.Cells.Item(19, 5).Value = 96411
$ws.Cells.Item(19, 6).Value = 4662260
$ws.Cells.Item(19, 7).Value = -349894
$ws.Cells.Item(19, 8).Value = -400985

# Row 20
$ws.Cells.Item(20, 4).Value = 13137817
$ws.Cells.Item(20, 5).Value = 26110199
$ws.Cells.Item(20, 6).Value = 38094752
$ws.Cells.Item(20, 7).Value = 7220430
$ws.Cells.Item(20, 8).Value = 11662673

# Row 21
$ws.Cells.Item(21, 4).Value = -1072342
$ws.Cells.Item(21, 5).Value = -1072342
$ws.Cells.Item(21, 6).Value = -3666060
$ws.Cells.Item(21, 7).Value = 0
$ws.Cells.Item(21, 8).Value = -662248

# Row 22
$ws.Cells.Item(22, 4).Value = 12065475
$ws.Cells.Item(22, 5).Value = 25037857
$ws.Cells.Item(22, 6).Value = 34428692
$ws.Cells.Item(22, 7).Value = 7220430
$ws.Cells.Item(22, 8).Value = 11000425

# Row 23
$ws.Cells.Item(23, 4).Value = 0
$ws.Cells.Item(23, 5).Value = 0
$ws.Cells.Item(23, 6).Value = 0
$ws.Cells.Item(23, 7).Value = 0
$ws.Cells.Item(23, 8).Value = 0

# Row 24
$ws.Cells.Item(24, 4).Value = 12065475
$ws.Cells.Item(24, 5).Value = 25037857
$ws.Cells.Item(24, 6).Value = 34428692
$ws.Cells.Item(24, 7).Value = 7220430
$ws.Cells.Item(24, 8).Value = 11000425

# Row 25
$ws.Cells.Item(25, 4).Value = 804
$ws.Cells.Item(25, 5).Value = 1669
$ws.Cells.Item(25, 6).Value = 2295
$ws.Cells.Item(25, 7).Value = 481
$ws.Cells.Item(25, 8).Value = 344

# Row 26
$ws.Cells.Item(26, 4).Value = 15000000
$ws.Cells.Item(26, 5).Value = 15000000
$ws.Cells.Item(26, 6).Value = 15000000
$ws.Cells.Item(26, 7).Value = 15000000
$ws.Cells.Item(26, 8).Value = 32000000

# Row 27
$ws.Cells.Item(27, 4).Value = 377
$ws.Cells.Item(27, 5).Value = 782
$ws.Cells.Item(27, 6).Value = 1076
$ws.Cells.Item(27, 7).Value = 226
$ws.Cells.Item(27, 8).Value = 344

